$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 updates
$ws.Range("G4").Value  = 1.22
$ws.Range("H4").Value  = 5.4
$ws.Range("I4").Value  = 9.25
$ws.Range("J4").Value  = 1.6
$ws.Range("K4").Value  = 2.72
$ws.Range("L4").Value  = 7.4
$ws.Range("U4").Value  = 1.7
$ws.Range("V4").Value  = 2.1
$ws.Range("Y4").Value  = 7.8
$ws.Range("AA4").Value = 8.25
$ws.Range("AB4").Value = 17.5
$ws.Range("AC4").Value = 20
$ws.Range("AD4").Value = 10.5
$ws.Range("AE4").Value = 16.5
$ws.Range("AF4").Value = 50
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 29
$ws.Range("AI4").Value = 65
$ws.Range("AJ4").Value = 24
$ws.Range("AK4").Value = 200
$ws.Range("AL4").Value = 80
$ws.Range("AN4").Value = 3.35
$ws.Range("AQ4").Value = 11.75
$ws.Range("AT4").Value = 3.95
$ws.Range("AU4").Value = 8
$ws.Range("AW4").Value = 10.5
$ws.Range("AY4").Value = 40
$ws.Range("AZ4").Value = 350
$ws.Range("BA4").Value = 300

# Row 5 update
$ws.Range("P5").Value = 4.2
